$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New header cells (row 1), columns M:P ---
$ws.Cells.Item(1, 13).Value = "drop_islands"
$ws.Cells.Item(1, 14).Value = "drop_ocean_watersheds"
$ws.Cells.Item(1, 15).Value = "drop_within_sea"
$ws.Cells.Item(1, 16).Value = "drop_low_flow"

# --- New data cells (rows 2:63), columns M:P ---
$ws.Cells.Item(2, 13).Value = 1
$ws.Cells.Item(2, 14).Value = 1
$ws.Cells.Item(2, 15).Value = 1
$ws.Cells.Item(2, 16).Value = 1
$ws.Cells.Item(3, 13).Value = 1
$ws.Cells.Item(3, 14).Value = 1
$ws.Cells.Item(3, 15).Value = 1
$ws.Cells.Item(3, 16).Value = 0
$ws.Cells.Item(4, 13).Value = 1
$ws.Cells.Item(4, 14).Value = 1
$ws.Cells.Item(4, 15).Value = 1
$ws.Cells.Item(4, 16).Value = 0
$ws.Cells.Item(5, 13).Value = 1
$ws.Cells.Item(5, 14).Value = 1
$ws.Cells.Item(5, 15).Value = 1
$ws.Cells.Item(5, 16).Value = 1
$ws.Cells.Item(6, 13).Value = 1
$ws.Cells.Item(6, 14).Value = 1
$ws.Cells.Item(6, 15).Value = 1
$ws.Cells.Item(6, 16).Value = 1
$ws.Cells.Item(7, 13).Value = 1
$ws.Cells.Item(7, 14).Value = 1
$ws.Cells.Item(7, 15).Value = 1
$ws.Cells.Item(7, 16).Value = 1
$ws.Cells.Item(8, 13).Value = 1
$ws.Cells.Item(8, 14).Value = 1
$ws.Cells.Item(8, 15).Value = 1
$ws.Cells.Item(8, 16).Value = 0
$ws.Cells.Item(9, 13).Value = 1
$ws.Cells.Item(9, 14).Value = 1
$ws.Cells.Item(9, 15).Value = 1
$ws.Cells.Item(9, 16).Value = 1
$ws.Cells.Item(10, 13).Value = 1
$ws.Cells.Item(10, 14).Value = 1
$ws.Cells.Item(10, 15).Value = 1
$ws.Cells.Item(10, 16).Value = 0
$ws.Cells.Item(11, 13).Value = 1
$ws.Cells.Item(11, 14).Value = 1
$ws.Cells.Item(11, 15).Value = 1
$ws.Cells.Item(11, 16).Value = 0
$ws.Cells.Item(12, 13).Value = 1
$ws.Cells.Item(12, 14).Value = 1
$ws.Cells.Item(12, 15).Value = 1
$ws.Cells.Item(12, 16).Value = 0
$ws.Cells.Item(13, 13).Value = 1
$ws.Cells.Item(13, 14).Value = 1
$ws.Cells.Item(13, 15).Value = 1
$ws.Cells.Item(13, 16).Value = 0
$ws.Cells.Item(14, 13).Value = 1
$ws.Cells.Item(14, 14).Value = 1
$ws.Cells.Item(14, 15).Value = 1
$ws.Cells.Item(14, 16).Value = 0
$ws.Cells.Item(15, 13).Value = 1
$ws.Cells.Item(15, 14).Value = 1
$ws.Cells.Item(15, 15).Value = 1
$ws.Cells.Item(15, 16).Value = 0
$ws.Cells.Item(16, 13).Value = 1
$ws.Cells.Item(16, 14).Value = 1
$ws.Cells.Item(16, 15).Value = 1
$ws.Cells.Item(16, 16).Value = 0
$ws.Cells.Item(17, 13).Value = 1
$ws.Cells.Item(17, 14).Value = 1
$ws.Cells.Item(17, 15).Value = 1
$ws.Cells.Item(17, 16).Value = 0
$ws.Cells.Item(18, 13).Value = 1
$ws.Cells.Item(18, 14).Value = 1
$ws.Cells.Item(18, 15).Value = 1
$ws.Cells.Item(18, 16).Value = 1
$ws.Cells.Item(19, 13).Value = 1
$ws.Cells.Item(19, 14).Value = 1
$ws.Cells.Item(19, 15).Value = 1
$ws.Cells.Item(19, 16).Value = 0
$ws.Cells.Item(20, 13).Value = 1
$ws.Cells.Item(20, 14).Value = 1
$ws.Cells.Item(20, 15).Value = 1
$ws.Cells.Item(20, 16).Value = 0
$ws.Cells.Item(21, 13).Value = 1
$ws.Cells.Item(21, 14).Value = 1
$ws.Cells.Item(21, 15).Value = 1
$ws.Cells.Item(21, 16).Value = 0
$ws.Cells.Item(22, 13).Value = 1
$ws.Cells.Item(22, 14).Value = 1
$ws.Cells.Item(22, 15).Value = 1
$ws.Cells.Item(22, 16).Value = 0
$ws.Cells.Item(23, 13).Value = 1
$ws.Cells.Item(23, 14).Value = 1
$ws.Cells.Item(23, 15).Value = 1
$ws.Cells.Item(23, 16).Value = 0
$ws.Cells.Item(24, 13).Value = 1
$ws.Cells.Item(24, 14).Value = 1
$ws.Cells.Item(24, 15).Value = 1
$ws.Cells.Item(24, 16).Value = 0
$ws.Cells.Item(25, 13).Value = 1
$ws.Cells.Item(25, 14).Value = 1
$ws.Cells.Item(25, 15).Value = 1
$ws.Cells.Item(25, 16).Value = 0
$ws.Cells.Item(26, 13).Value = 1
$ws.Cells.Item(26, 14).Value = 1
$ws.Cells.Item(26, 15).Value = 1
$ws.Cells.Item(26, 16).Value = 0
$ws.Cells.Item(27, 13).Value = 1
$ws.Cells.Item(27, 14).Value = 1
$ws.Cells.Item(27, 15).Value = 1
$ws.Cells.Item(27, 16).Value = 0
$ws.Cells.Item(28, 13).Value = 1
$ws.Cells.Item(28, 14).Value = 1
$ws.Cells.Item(28, 15).Value = 1
$ws.Cells.Item(28, 16).Value = 0
$ws.Cells.Item(29, 13).Value = 1
$ws.Cells.Item(29, 14).Value = 1
$ws.Cells.Item(29, 15).Value = 1
$ws.Cells.Item(29, 16).Value = 0
$ws.Cells.Item(30, 13).Value = 1
$ws.Cells.Item(30, 14).Value = 1
$ws.Cells.Item(30, 15).Value = 1
$ws.Cells.Item(30, 16).Value = 1
$ws.Cells.Item(31, 13).Value = 1
$ws.Cells.Item(31, 14).Value = 1
$ws.Cells.Item(31, 15).Value = 1
$ws.Cells.Item(31, 16).Value = 0
$ws.Cells.Item(32, 13).Value = 1
$ws.Cells.Item(32, 14).Value = 1
$ws.Cells.Item(32, 15).Value = 1
$ws.Cells.Item(32, 16).Value = 0
$ws.Cells.Item(33, 13).Value = 1
$ws.Cells.Item(33, 14).Value = 1
$ws.Cells.Item(33, 15).Value = 1
$ws.Cells.Item(33, 16).Value = 0
$ws.Cells.Item(34, 13).Value = 1
$ws.Cells.Item(34, 14).Value = 1
$ws.Cells.Item(34, 15).Value = 1
$ws.Cells.Item(34, 16).Value = 0
$ws.Cells.Item(35, 13).Value = 1
$ws.Cells.Item(35, 14).Value = 1
$ws.Cells.Item(35, 15).Value = 1
$ws.Cells.Item(35, 16).Value = 0
$ws.Cells.Item(36, 13).Value = 1
$ws.Cells.Item(36, 14).Value = 1
$ws.Cells.Item(36, 15).Value = 1
$ws.Cells.Item(36, 16).Value = 0
$ws.Cells.Item(37, 13).Value = 1
$ws.Cells.Item(37, 14).Value = 1
$ws.Cells.Item(37, 15).Value = 1
$ws.Cells.Item(37, 16).Value = 1
$ws.Cells.Item(38, 13).Value = 1
$ws.Cells.Item(38, 14).Value = 1
$ws.Cells.Item(38, 15).Value = 1
$ws.Cells.Item(38, 16).Value = 0
$ws.Cells.Item(39, 13).Value = 1
$ws.Cells.Item(39, 14).Value = 1
$ws.Cells.Item(39, 15).Value = 1
$ws.Cells.Item(39, 16).Value = 0
$ws.Cells.Item(40, 13).Value = 1
$ws.Cells.Item(40, 14).Value = 1
$ws.Cells.Item(40, 15).Value = 1
$ws.Cells.Item(40, 16).Value = 0
$ws.Cells.Item(41, 13).Value = 1
$ws.Cells.Item(41, 14).Value = 1
$ws.Cells.Item(41, 15).Value = 1
$ws.Cells.Item(41, 16).Value = 0
$ws.Cells.Item(42, 13).Value = 1
$ws.Cells.Item(42, 14).Value = 1
$ws.Cells.Item(42, 15).Value = 1
$ws.Cells.Item(42, 16).Value = 0
$ws.Cells.Item(43, 13).Value = 1
$ws.Cells.Item(43, 14).Value = 1
$ws.Cells.Item(43, 15).Value = 1
$ws.Cells.Item(43, 16).Value = 0
$ws.Cells.Item(44, 13).Value = 1
$ws.Cells.Item(44, 14).Value = 1
$ws.Cells.Item(44, 15).Value = 1
$ws.Cells.Item(44, 16).Value = 0
$ws.Cells.Item(45, 13).Value = 1
$ws.Cells.Item(45, 14).Value = 1
$ws.Cells.Item(45, 15).Value = 1
$ws.Cells.Item(45, 16).Value = 0
$ws.Cells.Item(46, 13).Value = 1
$ws.Cells.Item(46, 14).Value = 1
$ws.Cells.Item(46, 15).Value = 1
$ws.Cells.Item(46, 16).Value = 0
$ws.Cells.Item(47, 13).Value = 1
$ws.Cells.Item(47, 14).Value = 1
$ws.Cells.Item(47, 15).Value = 1
$ws.Cells.Item(47, 16).Value = 0
$ws.Cells.Item(48, 13).Value = 1
$ws.Cells.Item(48, 14).Value = 1
$ws.Cells.Item(48, 15).Value = 1
$ws.Cells.Item(48, 16).Value = 0
$ws.Cells.Item(49, 13).Value = 1
$ws.Cells.Item(49, 14).Value = 1
$ws.Cells.Item(49, 15).Value = 1
$ws.Cells.Item(49, 16).Value = 0
$ws.Cells.Item(50, 13).Value = 1
$ws.Cells.Item(50, 14).Value = 1
$ws.Cells.Item(50, 15).Value = 1
$ws.Cells.Item(50, 16).Value = 0
$ws.Cells.Item(51, 13).Value = 1
$ws.Cells.Item(51, 14).Value = 1
$ws.Cells.Item(51, 15).Value = 1
$ws.Cells.Item(51, 16).Value = 0
$ws.Cells.Item(52, 13).Value = 1
$ws.Cells.Item(52, 14).Value = 1
$ws.Cells.Item(52, 15).Value = 1
$ws.Cells.Item(52, 16).Value = 0
$ws.Cells.Item(53, 13).Value = 1
$ws.Cells.Item(53, 14).Value = 1
$ws.Cells.Item(53, 15).Value = 1
$ws.Cells.Item(53, 16).Value = 0
$ws.Cells.Item(54, 13).Value = 1
$ws.Cells.Item(54, 14).Value = 1
$ws.Cells.Item(54, 15).Value = 1
$ws.Cells.Item(54, 16).Value = 0
$ws.Cells.Item(55, 13).Value = 1
$ws.Cells.Item(55, 14).Value = 1
$ws.Cells.Item(55, 15).Value = 1
$ws.Cells.Item(55, 16).Value = 0
$ws.Cells.Item(56, 13).Value = 1
$ws.Cells.Item(56, 14).Value = 1
$ws.Cells.Item(56, 15).Value = 1
$ws.Cells.Item(56, 16).Value = 0
$ws.Cells.Item(57, 13).Value = 1
$ws.Cells.Item(57, 14).Value = 1
$ws.Cells.Item(57, 15).Value = 1
$ws.Cells.Item(57, 16).Value = 0
$ws.Cells.Item(58, 13).Value = 1
$ws.Cells.Item(58, 14).Value = 1
$ws.Cells.Item(58, 15).Value = 1
$ws.Cells.Item(58, 16).Value = 0
$ws.Cells.Item(59, 13).Value = 1
$ws.Cells.Item(59, 14).Value = 1
$ws.Cells.Item(59, 15).Value = 1
$ws.Cells.Item(59, 16).Value = 0
$ws.Cells.Item(60, 13).Value = 1
$ws.Cells.Item(60, 14).Value = 1
$ws.Cells.Item(60, 15).Value = 1
$ws.Cells.Item(60, 16).Value = 0
$ws.Cells.Item(61, 13).Value = 1
$ws.Cells.Item(61, 14).Value = 1
$ws.Cells.Item(61, 15).Value = 1
$ws.Cells.Item(61, 16).Value = 0
$ws.Cells.Item(62, 13).Value = 1
$ws.Cells.Item(62, 14).Value = 1
$ws.Cells.Item(62, 15).Value = 1
$ws.Cells.Item(62, 16).Value = 0
$ws.Cells.Item(63, 13).Value = 1
$ws.Cells.Item(63, 14).Value = 1
$ws.Cells.Item(63, 15).Value = 1
$ws.Cells.Item(63, 16).Value = 0

# --- Column widths (engine quantizes ColumnWidth to 1/6-character
#     increments, so these land within ~0.001 of the target bestFit widths) ---
$ws.Columns.Item(4).ColumnWidth = 9
$ws.Columns.Item(12).ColumnWidth = 12
$ws.Columns.Item(13).ColumnWidth = 10.6
$ws.Columns.Item(14).ColumnWidth = 20.6
$ws.Columns.Item(15).ColumnWidth = 14
$ws.Columns.Item(16).ColumnWidth = 12.6

# --- Restore frozen header pane, move viewport/selection to the top
#     of the sheet and select Q2 (the new first empty cell) ---
$ws.Range("A2").Select()
$win = $excel.ActiveWindow
$win.FreezePanes = $true
$ws.Range("Q2").Select()
